# Applies targeted odds updates to Sheet1 of the workbook as described in the
# commit diff (167 individual cell value updates across 16 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 20).Value = 9.5  # T2: 10 -> 9.5
$ws.Cells.Item(2, 21).Value = 13  # U2: 12 -> 13
$ws.Cells.Item(2, 23).Value = 23  # W2: 21 -> 23
$ws.Cells.Item(3, 7).Value = 2.4  # G3: 2.45 -> 2.4
$ws.Cells.Item(3, 9).Value = 3.5  # I3: 3.4 -> 3.5
$ws.Cells.Item(3, 11).Value = 4.75  # K3: 5 -> 4.75
$ws.Cells.Item(3, 31).Value = 6.5  # AE3: 6 -> 6.5
$ws.Cells.Item(4, 8).Value = 3.6  # H4: 3.75 -> 3.6
$ws.Cells.Item(4, 9).Value = 4.33  # I4: 4.2 -> 4.33
$ws.Cells.Item(4, 20).Value = 8.5  # T4: 8 -> 8.5
$ws.Cells.Item(4, 33).Value = 15  # AG4: 13 -> 15
$ws.Cells.Item(6, 7).Value = 1.45  # G6: 1.4 -> 1.45
$ws.Cells.Item(6, 8).Value = 4.1  # H6: 4.33 -> 4.1
$ws.Cells.Item(6, 9).Value = 7  # I6: 7.5 -> 7
$ws.Cells.Item(6, 22).Value = 9  # V6: 8.5 -> 9
$ws.Cells.Item(6, 26).Value = 15  # Z6: 17 -> 15
$ws.Cells.Item(6, 27).Value = 8.5  # AA6: 9 -> 8.5
$ws.Cells.Item(6, 33).Value = 21  # AG6: 23 -> 21
$ws.Cells.Item(9, 7).Value = 2.3  # G9: 2.27 -> 2.3
$ws.Cells.Item(9, 8).Value = 2.67  # H9: 2.7 -> 2.67
$ws.Cells.Item(9, 11).Value = 5.8  # K9: 5.9 -> 5.8
$ws.Cells.Item(9, 12).Value = 1.39  # L9: 1.38 -> 1.39
$ws.Cells.Item(9, 13).Value = 2.77  # M9: 2.8 -> 2.77
$ws.Cells.Item(9, 14).Value = 2.15  # N9: 2.12 -> 2.15
$ws.Cells.Item(9, 15).Value = 1.62  # O9: 1.65 -> 1.62
$ws.Cells.Item(9, 16).Value = 1.45  # P9: 1.44 -> 1.45
$ws.Cells.Item(9, 17).Value = 2.55  # Q9: 2.6 -> 2.55
$ws.Cells.Item(9, 26).Value = 5.8  # Z9: 5.9 -> 5.8
$ws.Cells.Item(9, 27).Value = 5.2  # AA9: 5.3 -> 5.2
$ws.Cells.Item(9, 31).Value = 9.5  # AE9: 9.75 -> 9.5
$ws.Cells.Item(10, 8).Value = 3.3  # H10: 3.35 -> 3.3
$ws.Cells.Item(10, 9).Value = 4.1  # I10: 4.05 -> 4.1
$ws.Cells.Item(10, 12).Value = 1.38  # L10: 1.37 -> 1.38
$ws.Cells.Item(10, 13).Value = 2.82  # M10: 2.85 -> 2.82
$ws.Cells.Item(10, 14).Value = 2.12  # N10: 2.1 -> 2.12
$ws.Cells.Item(10, 18).Value = 1.95  # R10: 1.93 -> 1.95
$ws.Cells.Item(10, 19).Value = 1.75  # S10: 1.78 -> 1.75
$ws.Cells.Item(10, 20).Value = 6.1  # T10: 6.2 -> 6.1
$ws.Cells.Item(10, 24).Value = 16.5  # X10: 16 -> 16.5
$ws.Cells.Item(10, 30).Value = 900  # AD10: 800 -> 900
$ws.Cells.Item(11, 7).Value = 1.83  # G11: 1.85 -> 1.83
$ws.Cells.Item(11, 9).Value = 4.33  # I11: 4.2 -> 4.33
$ws.Cells.Item(11, 12).Value = 1.33  # L11: 1.36 -> 1.33
$ws.Cells.Item(11, 13).Value = 3.25  # M11: 3 -> 3.25
$ws.Cells.Item(11, 18).Value = 2  # R11: 1.91 -> 2
$ws.Cells.Item(11, 19).Value = 1.73  # S11: 1.8 -> 1.73
$ws.Cells.Item(11, 24).Value = 15  # X11: 17 -> 15
$ws.Cells.Item(11, 31).Value = 11  # AE11: 10 -> 11
$ws.Cells.Item(11, 34).Value = 51  # AH11: 41 -> 51
$ws.Cells.Item(12, 7).Value = 3.4  # G12: 3.25 -> 3.4
$ws.Cells.Item(12, 8).Value = 3.7  # H12: 3.6 -> 3.7
$ws.Cells.Item(12, 9).Value = 2  # I12: 2.05 -> 2
$ws.Cells.Item(12, 11).Value = 9.5  # K12: 10 -> 9.5
$ws.Cells.Item(12, 20).Value = 9  # T12: 8.5 -> 9
$ws.Cells.Item(12, 21).Value = 17  # U12: 15 -> 17
$ws.Cells.Item(12, 26).Value = 9.5  # Z12: 9 -> 9.5
$ws.Cells.Item(12, 28).Value = 19  # AB12: 17 -> 19
$ws.Cells.Item(12, 32).Value = 9  # AF12: 9.5 -> 9
$ws.Cells.Item(12, 34).Value = 17  # AH12: 19 -> 17
$ws.Cells.Item(12, 35).Value = 17  # AI12: 19 -> 17
$ws.Cells.Item(12, 36).Value = 29  # AJ12: 34 -> 29
$ws.Cells.Item(17, 7).Value = 3.35  # G17: 3.75 -> 3.35
$ws.Cells.Item(17, 8).Value = 3.1  # H17: 3.2 -> 3.1
$ws.Cells.Item(17, 9).Value = 2.12  # I17: 1.95 -> 2.12
$ws.Cells.Item(17, 12).Value = 1.37  # L17: 1.36 -> 1.37
$ws.Cells.Item(17, 13).Value = 2.62  # M17: 2.65 -> 2.62
$ws.Cells.Item(17, 14).Value = 2.07  # N17: 2.05 -> 2.07
$ws.Cells.Item(17, 16).Value = 1.44  # P17: 1.42 -> 1.44
$ws.Cells.Item(17, 17).Value = 2.42  # Q17: 2.47 -> 2.42
$ws.Cells.Item(17, 18).Value = 1.85  # R17: 1.87 -> 1.85
$ws.Cells.Item(17, 20).Value = 8.5  # T17: 9.5 -> 8.5
$ws.Cells.Item(17, 21).Value = 16.5  # U17: 19.5 -> 16.5
$ws.Cells.Item(17, 22).Value = 12  # V17: 13 -> 12
$ws.Cells.Item(17, 23).Value = 45  # W17: 60 -> 45
$ws.Cells.Item(17, 24).Value = 35  # X17: 40 -> 35
$ws.Cells.Item(17, 25).Value = 45  # Y17: 50 -> 45
$ws.Cells.Item(17, 26).Value = 7.9  # Z17: 8 -> 7.9
$ws.Cells.Item(17, 27).Value = 6.1  # AA17: 6.3 -> 6.1
$ws.Cells.Item(17, 28).Value = 15.5  # AB17: 16.5 -> 15.5
$ws.Cells.Item(17, 29).Value = 80  # AC17: 90 -> 80
$ws.Cells.Item(17, 31).Value = 6.6  # AE17: 6.3 -> 6.6
$ws.Cells.Item(17, 32).Value = 9.75  # AF17: 8.5 -> 9.75
$ws.Cells.Item(17, 33).Value = 9  # AG17: 8.5 -> 9
$ws.Cells.Item(17, 34).Value = 20  # AH17: 17 -> 20
$ws.Cells.Item(17, 35).Value = 18.5  # AI17: 17 -> 18.5
$ws.Cells.Item(19, 8).Value = 2.92  # H19: 2.95 -> 2.92
$ws.Cells.Item(19, 12).Value = 1.36  # L19: 1.35 -> 1.36
$ws.Cells.Item(19, 13).Value = 2.65  # M19: 2.7 -> 2.65
$ws.Cells.Item(19, 14).Value = 2.05  # N19: 2.02 -> 2.05
$ws.Cells.Item(19, 15).Value = 1.6  # O19: 1.62 -> 1.6
$ws.Cells.Item(19, 16).Value = 1.5  # P19: 1.47 -> 1.5
$ws.Cells.Item(19, 17).Value = 2.27  # Q19: 2.32 -> 2.27
$ws.Cells.Item(19, 18).Value = 1.78  # R19: 1.75 -> 1.78
$ws.Cells.Item(19, 19).Value = 1.83  # S19: 1.85 -> 1.83
$ws.Cells.Item(19, 20).Value = 9  # T19: 9.5 -> 9
$ws.Cells.Item(19, 21).Value = 18  # U19: 18.5 -> 18
$ws.Cells.Item(19, 22).Value = 11.75  # V19: 11.5 -> 11.75
$ws.Cells.Item(19, 24).Value = 35  # X19: 32 -> 35
$ws.Cells.Item(19, 26).Value = 7.6  # Z19: 7.8 -> 7.6
$ws.Cells.Item(19, 31).Value = 6.8  # AE19: 6.7 -> 6.8
$ws.Cells.Item(23, 7).Value = 3  # G23: 2.77 -> 3
$ws.Cells.Item(23, 8).Value = 3.35  # H23: 3.25 -> 3.35
$ws.Cells.Item(23, 9).Value = 2.07  # I23: 2.27 -> 2.07
$ws.Cells.Item(23, 14).Value = 1.8  # N23: 1.82 -> 1.8
$ws.Cells.Item(23, 15).Value = 1.8  # O23: 1.78 -> 1.8
$ws.Cells.Item(23, 20).Value = 8.5  # T23: 7.8 -> 8.5
$ws.Cells.Item(23, 21).Value = 13  # U23: 12 -> 13
$ws.Cells.Item(23, 22).Value = 9.25  # V23: 8.75 -> 9.25
$ws.Cells.Item(23, 23).Value = 29  # W23: 26 -> 29
$ws.Cells.Item(23, 24).Value = 20  # X23: 18.5 -> 20
$ws.Cells.Item(23, 25).Value = 25  # Y23: 24 -> 25
$ws.Cells.Item(23, 26).Value = 10.25  # Z23: 10 -> 10.25
$ws.Cells.Item(23, 27).Value = 5.8  # AA23: 5.6 -> 5.8
$ws.Cells.Item(23, 28).Value = 11.5  # AB23: 11 -> 11.5
$ws.Cells.Item(23, 31).Value = 6.7  # AE23: 7.1 -> 6.7
$ws.Cells.Item(23, 32).Value = 8.75  # AF23: 9.75 -> 8.75
$ws.Cells.Item(23, 33).Value = 7.5  # AG23: 7.7 -> 7.5
$ws.Cells.Item(23, 34).Value = 15.5  # AH23: 18.5 -> 15.5
$ws.Cells.Item(23, 35).Value = 13.5  # AI23: 14.5 -> 13.5
$ws.Cells.Item(25, 7).Value = 2.05  # G25: 2.15 -> 2.05
$ws.Cells.Item(25, 9).Value = 3.7  # I25: 3.5 -> 3.7
$ws.Cells.Item(25, 10).Value = 1.1  # J25: 1.11 -> 1.1
$ws.Cells.Item(25, 11).Value = 7  # K25: 6.5 -> 7
$ws.Cells.Item(25, 18).Value = 2.2  # R25: 2.1 -> 2.2
$ws.Cells.Item(25, 19).Value = 1.62  # S25: 1.67 -> 1.62
$ws.Cells.Item(25, 21).Value = 8.5  # U25: 9 -> 8.5
$ws.Cells.Item(25, 22).Value = 9.5  # V25: 10 -> 9.5
$ws.Cells.Item(25, 33).Value = 15  # AG25: 13 -> 15
$ws.Cells.Item(25, 35).Value = 41  # AI25: 34 -> 41
$ws.Cells.Item(25, 36).Value = 51  # AJ25: 41 -> 51
$ws.Cells.Item(26, 7).Value = 3.1  # G26: 3 -> 3.1
$ws.Cells.Item(26, 9).Value = 2.35  # I26: 2.4 -> 2.35
$ws.Cells.Item(26, 22).Value = 12  # V26: 11 -> 12
$ws.Cells.Item(26, 23).Value = 34  # W26: 29 -> 34
$ws.Cells.Item(26, 33).Value = 9.5  # AG26: 10 -> 9.5
$ws.Cells.Item(26, 34).Value = 21  # AH26: 23 -> 21
$ws.Cells.Item(26, 36).Value = 29  # AJ26: 34 -> 29
$ws.Cells.Item(33, 7).Value = 1.6  # G33: 1.62 -> 1.6
$ws.Cells.Item(33, 9).Value = 4.33  # I33: 4.2 -> 4.33
$ws.Cells.Item(33, 10).Value = 1.01  # J33: 1.02 -> 1.01
$ws.Cells.Item(33, 11).Value = 13  # K33: 12 -> 13
$ws.Cells.Item(33, 12).Value = 1.1  # L33: 1.11 -> 1.1
$ws.Cells.Item(33, 13).Value = 6.5  # M33: 6 -> 6.5
$ws.Cells.Item(33, 14).Value = 1.36  # N33: 1.4 -> 1.36
$ws.Cells.Item(33, 15).Value = 3  # O33: 2.75 -> 3
$ws.Cells.Item(33, 25).Value = 17  # Y33: 19 -> 17
$ws.Cells.Item(33, 26).Value = 23  # Z33: 21 -> 23
$ws.Cells.Item(33, 34).Value = 51  # AH33: 41 -> 51
$ws.Cells.Item(39, 11).Value = 9.5  # K39: 10 -> 9.5
$ws.Cells.Item(39, 18).Value = 1.8  # R39: 1.75 -> 1.8
$ws.Cells.Item(39, 19).Value = 1.95  # S39: 2 -> 1.95
$ws.Cells.Item(39, 20).Value = 8  # T39: 8.5 -> 8
$ws.Cells.Item(39, 26).Value = 9.5  # Z39: 10 -> 9.5
$ws.Cells.Item(39, 28).Value = 15  # AB39: 13 -> 15
$ws.Cells.Item(39, 31).Value = 8.5  # AE39: 9 -> 8.5
$ws.Cells.Item(39, 33).Value = 11  # AG39: 10 -> 11
$ws.Cells.Item(39, 36).Value = 34  # AJ39: 29 -> 34
$ws.Cells.Item(41, 7).Value = 1.36  # G41: 1.33 -> 1.36
$ws.Cells.Item(41, 8).Value = 4.75  # H41: 5 -> 4.75
$ws.Cells.Item(41, 9).Value = 8  # I41: 8.5 -> 8
$ws.Cells.Item(41, 18).Value = 1.95  # R41: 2 -> 1.95
$ws.Cells.Item(41, 19).Value = 1.8  # S41: 1.75 -> 1.8
$ws.Cells.Item(41, 20).Value = 7.5  # T41: 7 -> 7.5
$ws.Cells.Item(41, 23).Value = 9  # W41: 8.5 -> 9
$ws.Cells.Item(41, 27).Value = 9.5  # AA41: 10 -> 9.5
$ws.Cells.Item(41, 29).Value = 51  # AC41: 67 -> 51
$ws.Cells.Item(41, 30).Value = 351  # AD41: 401 -> 351

$wb.Save()
